$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" sheet: A1 conversion text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.92 = 6954.42 pesos`n✅ 6954.42 pesos = 1.91 = 966.59 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet: N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 521.97
$ws2.Range("O10").Value = 3630
$ws2.Range("N12").Value = 3647.77
$ws2.Range("O12").Value = 507
